# Apply crypto price/volume updates to Sheet1 (cryptos.xlsx)
# Generated from the commit diff: rows 2-51, columns D (Price) and E (Volume(1h)).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D values that read as plain decimals (e.g. "556.27") get a leading
# apostrophe so Excel keeps them as text, matching the source data's string type
# (values like "59.279.65" already have two dots and are never auto-converted).

$ws.Range("D2").Value = '59.279.65'
$ws.Range("E2").Value = '  -1.54%  '

$ws.Range("D3").Value = '2.573.98'
$ws.Range("E3").Value = '  -1.85%  '

$ws.Range("E4").Value = '  -0.09%  '

$ws.Range("D5").Value = '''556.27'
$ws.Range("E5").Value = '  -2.04%  '

$ws.Range("D6").Value = '''141.51'
$ws.Range("E6").Value = '  -2.61%  '

$ws.Range("E7").Value = '  +0.23%  '

$ws.Range("E8").Value = '  -1.88%  '

$ws.Range("D9").Value = '2.579.15'
$ws.Range("E9").Value = '  -2.41%  '

$ws.Range("D10").Value = '''6.76'
$ws.Range("E10").Value = '  -1.04%  '

$ws.Range("E11").Value = '  -0.72%  '

$ws.Range("D12").Value = '''0.165'
$ws.Range("E12").Value = '  +11.30%  '

$ws.Range("D13").Value = '''0.351'

$ws.Range("D14").Value = '3.028.89'

$ws.Range("D15").Value = '59.270.42'
$ws.Range("E15").Value = '  -1.57%  '

$ws.Range("D16").Value = '''22.91'
$ws.Range("E16").Value = '  +4.15%  '

$ws.Range("E17").Value = '  -0.41%  '

$ws.Range("D18").Value = '2.573.25'
$ws.Range("E18").Value = '  -2.45%  '

$ws.Range("D19").Value = '''4.54'
$ws.Range("E19").Value = '  +0.37%  '

$ws.Range("D20").Value = '''337.96'
$ws.Range("E20").Value = '  -0.98%  '

$ws.Range("E21").Value = '  -0.76%  '

$ws.Range("E22").Value = '  +1.30%  '

$ws.Range("E23").Value = '  -0.11%  '

$ws.Range("E24").Value = '  +9.34%  '

$ws.Range("D25").Value = '''62.66'
$ws.Range("E25").Value = '  -4.99%  '

$ws.Range("E26").Value = '  +0.10%  '

$ws.Range("E27").Value = '  -2.53%  '

$ws.Range("D28").Value = '''7.39'
$ws.Range("E28").Value = '  +0.17%  '

$ws.Range("E29").Value = '  -3.65%  '

$ws.Range("E30").Value = '  +0.07%  '

$ws.Range("D31").Value = '''6.19'
$ws.Range("E31").Value = '  +1.12%  '

$ws.Range("E32").Value = '  -1.81%  '

$ws.Range("D33").Value = '''159.17'
$ws.Range("E33").Value = '  +0.61%  '

$ws.Range("D34").Value = '''19.06'
$ws.Range("E34").Value = '  -0.32%  '

$ws.Range("E35").Value = '  -0.41%  '

$ws.Range("E36").Value = '  +1.30%  '

$ws.Range("D37").Value = '''0.894'

$ws.Range("E38").Value = '  -0.30%  '

$ws.Range("D39").Value = '''0.852'
$ws.Range("E39").Value = '  -3.37%  '

$ws.Range("E40").Value = '  -2.24%  '

$ws.Range("E41").Value = '  +1.08%  '

$ws.Range("D42").Value = '''290.18'
$ws.Range("E42").Value = '  -2.91%  '

$ws.Range("D43").Value = '''138.74'
$ws.Range("E43").Value = '  +8.78%  '

$ws.Range("E44").Value = '  +0.38%  '

$ws.Range("E45").Value = '  -1.37%  '

$ws.Range("E46").Value = '  -1.59%  '

$ws.Range("E47").Value = '  -0.12%  '

$ws.Range("D48").Value = '''0.0530'
$ws.Range("E48").Value = '  -2.41%  '

$ws.Range("E49").Value = '  -0.42%  '

$ws.Range("E50").Value = '  -0.03%  '

$ws.Range("D51").Value = '1.936.94'
$ws.Range("E51").Value = '  -0.99%  '
